$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF")
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the style of the existing header cell (H1: bold, bordered, centered)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Fill in the data for columns I (I0) and J (IF) for rows 2-15
$data = @(
    @(8, 8),
    @(5, 6),
    @(9, 9),
    @(8, 8),
    @(8, 9),
    @(6, 6),
    @(8, 9),
    @(8, 8),
    @(8, 8),
    @(1, 9),
    @(6, 8),
    @(1, 3),
    @(1, 2),
    @(6, 6)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
